$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Date placeholder: "16/06/2021" -> "20/06/2021" on the slide master
#    and on every slide layout (ppPlaceholderDate = 16).
# ---------------------------------------------------------------------
function Update-DatePlaceholder($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $shp = $container.Shapes.Item($i)
        if ($shp.Type -eq 14) {
            if ($shp.PlaceholderFormat.Type -eq 16) {
                if ($shp.HasTextFrame -eq -1) {
                    $shp.TextFrame.TextRange.Text = "20/06/2021"
                }
            }
        }
    }
}

$master = $p.SlideMaster
Update-DatePlaceholder $master

for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    Update-DatePlaceholder $layout
}

# ---------------------------------------------------------------------
# 2) Slide 19 ("Loi de Biot :"), shape "ZoneTexte 5": reword + resize.
# ---------------------------------------------------------------------
$s19 = $p.Slides.Item(19)
$shp = $s19.Shapes.Item(4)

$tr = $shp.TextFrame.TextRange
$para1 = $tr.Paragraphs(1)

# Reset to an unrelated placeholder first so the runtime's prefix-diff
# run-splitter does not fragment the first (identically formatted) run.
$para1.Text = "PLACEHOLDER_TEXT_XYZ"
$para1.Text = "Valeur tabulée :  [α] =12.5 °."
[void]$para1.InsertAfter("mL")
[void]$para1.InsertAfter("/g/dm pour l’acide (+) ")
[void]$para1.InsertAfter("tratrique")

$emuPerPt = 12700
$shp.Left = 979889 / $emuPerPt
$shp.Width = 10015556 / $emuPerPt
$shp.Height = 830997 / $emuPerPt
